# Auto-generated: update cryptos worksheet Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D) updates ---
$ws.Range("D2").Value = "27.302.12"
$ws.Range("D3").Value = "1.653.05"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D5").Value = "'218.98"
$ws.Range("D6").Value = "'0.510"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D10").Value = "'20.27"
$ws.Range("D11").Value = "'0.0847"
$ws.Range("D13").Value = "1.643.70"
$ws.Range("D15").Value = "'0.543"
$ws.Range("D16").Value = "'67.99"
$ws.Range("D17").Value = "27.269.30"
$ws.Range("D19").Value = "'221.86"
$ws.Range("D20").Value = "'0.999"
$ws.Range("D22").Value = "'4.46"
$ws.Range("D23").Value = "'2.50"
$ws.Range("D24").Value = "'9.29"
$ws.Range("D25").Value = "'147.98"
$ws.Range("D26").Value = "'1.00"
$ws.Range("D29").Value = "'15.90"
$ws.Range("D32").Value = "'3.37"
$ws.Range("D35").Value = "1.269.45"
$ws.Range("D38").Value = "'0.544"
$ws.Range("D39").Value = "'0.846"
$ws.Range("D41").Value = "'0.813"
$ws.Range("D42").Value = "'5.40"
$ws.Range("D44").Value = "1.792.48"
$ws.Range("D45").Value = "'63.37"
$ws.Range("D46").Value = "'92.64"
$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("D49").Value = "'0.0515"
$ws.Range("D50").Value = "'7.70"
$ws.Range("D51").Value = "'0.0979"

# --- Volume(1h) column (E) updates ---
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("E6").Value = "  +2.38%  "
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("E8").Value = "  +1.70%  "
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("E10").Value = "  +4.63%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("E19").Value = "  +1.98%  "
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("E21").Value = "  +2.36%  "
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("E23").Value = "  +2.04%  "
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("E37").Value = "  +3.32%  "
$ws.Range("E38").Value = "  +2.27%  "
$ws.Range("E39").Value = "  +2.51%  "
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("E43").Value = "  +6.62%  "
$ws.Range("E45").Value = "  +2.65%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("E48").Value = "  +13.08%  "
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("E51").Value = "  +0.85%  "

